$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.675.26'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '1.965.09'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '244.85'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('D7').Value = '59.40'
$ws.Range('E7').Value = '  +2.47%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('D10').Value = '0.0815'
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '22.34'
$ws.Range('E12').Value = '  +4.12%  '
$ws.Range('D13').Value = '2.253.20'
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').Value = '0.831'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '13.74'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('E16').Value = '  +1.08%  '
$ws.Range('D17').Value = '1.964.97'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').Value = '36.582.55'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').Value = '70.06'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').Value = '0.0₃0860'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '229.42'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').Value = '5.08'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '2.46'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +3.72%  '
$ws.Range('D26').Value = '9.24'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('E27').Value = '  +12.01%  '
$ws.Range('D28').Value = '160.20'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').Value = '19.41'
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('E32').Value = '  +1.44%  '
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('D34').Value = '4.29'
$ws.Range('E34').Value = '  +1.24%  '
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = '2.26'
$ws.Range('E36').Value = '  +6.26%  '
$ws.Range('D37').Value = '3.42'
$ws.Range('E37').Value = '  +12.98%  '
$ws.Range('E38').Value = '  -2.91%  '
$ws.Range('D39').Value = '1.78'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').Value = '0.0985'
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D44').Value = '16.14'
$ws.Range('E44').Value = '  +0.96%  '
$ws.Range('D45').Value = '1.359.92'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('D47').Value = '87.96'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('D50').Value = '2.143.35'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('D51').Value = '43.90'
$ws.Range('E51').Value = '  -3.09%  '
